$d = $word.ActiveDocument
$r = $d.Paragraphs.Item(1).Range
$r.Collapse(1)
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512"><pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/></Relationships></pkg:xmlData></pkg:part><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr><w:t>M E M O R A N DU M</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr><w:t>31  DE</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr><w:t xml:space="preserve"> DICIEMBRE 2021</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr><w:t xml:space="preserve">PARA:     </w:t></w:r><w:r><w:rPr><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr><w:t xml:space="preserve">CENTRAL </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr><w:t xml:space="preserve">DE: </w:t></w:r><w:r><w:rPr><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr><w:tab/><w:t xml:space="preserve">   DIRECCION GENERAL </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr><w:t>ASUNTO</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="44"/><w:szCs w:val="44"/><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve">:  F A L T A S  </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:rPr><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:rPr><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr><w:tab/><w:t xml:space="preserve"> SE LES RECUERDA QUE AL ACUMULAR</w:t></w:r><w:r><w:rPr><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr><w:t xml:space="preserve">                </w:t></w:r><w:r><w:rPr><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="44"/><w:szCs w:val="44"/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve">4 </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b/><w:sz w:val="44"/><w:szCs w:val="44"/><w:u w:val="single"/></w:rPr><w:t>FALT</w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="44"/><w:szCs w:val="44"/><w:u w:val="single"/></w:rPr><w:t>AS  EN</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:b/><w:sz w:val="44"/><w:szCs w:val="44"/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve"> UN </w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="44"/><w:szCs w:val="44"/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve">MES </w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="44"/><w:szCs w:val="44"/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="44"/><w:szCs w:val="44"/><w:u w:val="single"/></w:rPr><w:t>SIN JUSTIFICAR</w:t></w:r><w:r><w:rPr><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr><w:t xml:space="preserve"> ES BAJA AUTOMATICAMENTE </w:t></w:r></w:p><w:p><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr></w:pPr></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r.InsertXML($xml)
